$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (the data row) into a new row 3, carrying over values,
# shared-string references and cell styles/number formats.
$ws.Rows("2:2").Copy()
$ws.Rows("3:3").Insert()

# Move the view back to the top-left (A1) and place the active selection
# on A7, matching the saved worksheet view state.
$ws.Range("A7").Select()
